# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Column D cells whose new text is a plain decimal number (e.g. "0.999") are
# written with a leading apostrophe so Excel keeps them as literal text
# (matching the sheet's original inline-string cells) instead of silently
# converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.506.77'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '3.444.28'
$ws.Range("E3").Value = '  +1.92%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''580.13'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("D6").Value = '''149.13'
$ws.Range("E6").Value = '  +9.26%  '
$ws.Range("D7").Value = '3.444.52'
$ws.Range("E7").Value = '  +1.94%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '''0.475'
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("E10").Value = '  +3.04%  '
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("D12").Value = '''0.392'
$ws.Range("E12").Value = '  +0.95%  '
$ws.Range("D13").Value = '4.028.27'
$ws.Range("E13").Value = '  +1.74%  '
$ws.Range("D14").Value = '''28.09'
$ws.Range("E14").Value = '  +6.81%  '
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("D17").Value = '3.441.39'
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").Value = '61.551.16'
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").Value = '''6.33'
$ws.Range("E19").Value = '  +8.64%  '
$ws.Range("D20").Value = '''14.36'
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("D21").Value = '''9.49'
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").Value = '''390.37'
$ws.Range("E22").Value = '  +3.63%  '
$ws.Range("D23").Value = '''0.571'
$ws.Range("E23").Value = '  +2.65%  '
$ws.Range("D24").Value = '3.589.92'
$ws.Range("D25").Value = '''72.76'
$ws.Range("E25").Value = '  +1.95%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("E28").Value = '  -1.72%  '
$ws.Range("E29").Value = '  +6.65%  '
$ws.Range("D30").Value = '''7.82'
$ws.Range("E30").Value = '  +3.64%  '
$ws.Range("E32").Value = '  -13.60%  '
$ws.Range("D33").Value = '''8.27'
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("D34").Value = '''2.18'
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = '''24.01'
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").Value = '''5.30'
$ws.Range("E37").Value = '  +1.76%  '
$ws.Range("D38").Value = '''7.09'
$ws.Range("E38").Value = '  +3.37%  '
$ws.Range("D39").Value = '''1.57'
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("D40").Value = '''165.87'
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("D41").Value = '''0.0794'
$ws.Range("E41").Value = '  +4.97%  '
$ws.Range("D42").Value = '''26.30'
$ws.Range("E42").Value = '  +8.99%  '
$ws.Range("D43").Value = '''0.795'
$ws.Range("E43").Value = '  +3.06%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '''4.50'
$ws.Range("E44").Value = '  +1.89%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''1.00'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").Value = '''42.28'
$ws.Range("E46").Value = '  +1.77%  '
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").Value = '2.617.36'
$ws.Range("E48").Value = '  +7.78%  '
$ws.Range("E49").Value = '  -2.45%  '
$ws.Range("D50").Value = '''7.10'
$ws.Range("E50").Value = '  +4.30%  '
$ws.Range("D51").Value = '''23.18'
$ws.Range("E51").Value = '  -0.86%  '
